$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")

# The "Scale categories" for the Predation/CrBr/CrBuCt scale columns in row 2-4
# (column J, "Units") no longer apply - clear those cells, keeping their style.
$ws.Range("J2:J4").ClearContents()

# Bring the sheet to the front and move the active selection from the far
# right (O1) onto the newly-cleared Units column so J2:J4 is highlighted,
# with the view scrolled back so column B is visible at the left edge.
$ws.Activate() | Out-Null
$ws.Range("B1").Select() | Out-Null
$ws.Range("J2:J4").Select() | Out-Null
